# Updated frequency-table run for publication.
# Rows 2-5 (bases A/C/G/T) x columns B:X (positions 1-23) get refreshed
# frequency values from the new pipeline run. Cells not listed below
# already read 0 and are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.956765412329864
$ws.Range("C2").Value = 0.0368294635708567
$ws.Range("D2").Value = 0.00480384307445957
$ws.Range("E2").Value = 0.00480384307445957
$ws.Range("F2").Value = 0.00240192153722978
$ws.Range("G2").Value = 0.00160128102481986
$ws.Range("I2").Value = 0.00720576461168935
$ws.Range("J2").Value = 0.00160128102481986
$ws.Range("K2").Value = 0.00240192153722978
$ws.Range("L2").Value = 0.00640512409927942
$ws.Range("M2").Value = 0.00160128102481986
$ws.Range("O2").Value = 0.00320256204963971
$ws.Range("P2").Value = 0.00400320256204964
$ws.Range("Q2").Value = 0.000800640512409928
$ws.Range("R2").Value = 0.00240192153722978
$ws.Range("S2").Value = 0.945556445156125
$ws.Range("T2").Value = 0.000800640512409928
$ws.Range("U2").Value = 0.0176140912730184
$ws.Range("V2").Value = 0.0200160128102482
$ws.Range("W2").Value = 0.00480384307445957
$ws.Range("X2").Value = 0.00640512409927942
# Row 3
$ws.Range("B3").Value = 0.032826261008807
$ws.Range("C3").Value = 0.00800640512409928
$ws.Range("D3").Value = 0.956765412329864
$ws.Range("E3").Value = 0.129703763010408
$ws.Range("F3").Value = 0.00960768614891913
$ws.Range("G3").Value = 0.964771817453963
$ws.Range("H3").Value = 0.0312249799839872
$ws.Range("I3").Value = 0.056044835868695
$ws.Range("J3").Value = 0.00160128102481986
$ws.Range("K3").Value = 0.204963971176942
$ws.Range("L3").Value = 0.0744595676541233
$ws.Range("M3").Value = 0.000800640512409928
$ws.Range("N3").Value = 0.99759807846277
$ws.Range("O3").Value = 0.926341072858287
$ws.Range("P3").Value = 0.00240192153722978
$ws.Range("R3").Value = 0.0056044835868695
$ws.Range("S3").Value = 0.00240192153722978
$ws.Range("T3").Value = 0.992794235388311
$ws.Range("U3").Value = 0.00320256204963971
$ws.Range("V3").Value = 0.000800640512409928
$ws.Range("W3").Value = 0.00240192153722978
# Row 4
$ws.Range("B4").Value = 0.0056044835868695
$ws.Range("C4").Value = 0.951961569255404
$ws.Range("D4").Value = 0.0368294635708567
$ws.Range("E4").Value = 0.00160128102481986
$ws.Range("F4").Value = 0.958366693354684
$ws.Range("G4").Value = 0.0312249799839872
$ws.Range("I4").Value = 0.00320256204963971
$ws.Range("J4").Value = 0.00320256204963971
$ws.Range("K4").Value = 0.0312249799839872
$ws.Range("L4").Value = 0.000800640512409928
$ws.Range("M4").Value = 0.99519615692554
$ws.Range("N4").Value = 0.00240192153722978
$ws.Range("O4").Value = 0.000800640512409928
$ws.Range("Q4").Value = 0.99919935948759
$ws.Range("R4").Value = 0.000800640512409928
$ws.Range("S4").Value = 0.00400320256204964
$ws.Range("U4").Value = 0.032826261008807
$ws.Range("V4").Value = 0.975980784627702
$ws.Range("W4").Value = 0.991993594875901
$ws.Range("X4").Value = 0.990392313851081
# Row 5
$ws.Range("B5").Value = 0.00480384307445957
$ws.Range("C5").Value = 0.00240192153722978
$ws.Range("D5").Value = 0.00160128102481986
$ws.Range("E5").Value = 0.863891112890312
$ws.Range("F5").Value = 0.0296236989591673
$ws.Range("G5").Value = 0.00240192153722978
$ws.Range("H5").Value = 0.968775020016013
$ws.Range("I5").Value = 0.933546837469976
$ws.Range("J5").Value = 0.993594875900721
$ws.Range("K5").Value = 0.760608486789431
$ws.Range("L5").Value = 0.917534027221777
$ws.Range("M5").Value = 0.00240192153722978
$ws.Range("O5").Value = 0.0696557245796637
$ws.Range("P5").Value = 0.993594875900721
$ws.Range("R5").Value = 0.991192954363491
$ws.Range("S5").Value = 0.0480384307445957
$ws.Range("T5").Value = 0.00640512409927942
$ws.Range("U5").Value = 0.946357085668535
$ws.Range("V5").Value = 0.00240192153722978
$ws.Range("W5").Value = 0.000800640512409928
$ws.Range("X5").Value = 0.00240192153722978
